$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# Version bump 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date bump
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank, now "Alvearie Team"
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 used to duplicate the "Contact" row; it becomes the Jurisdiction row
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was the second (duplicate) "Contact" row - remove it entirely, shifting
# everything below up by one row.
$meta.Rows(11).Delete()

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")

# Root Extension row: Short/Definition now use the profile-specific title & description
# instead of the generic "Extension" / "An Extension" placeholders.
$elements.Range("K2").Value = "Psychiatric Coverage Indicator"
$elements.Range("L2").Value = "Indicates whether the member has mental health and substance abuse benefit coverage: Y or N"
